{"js": "const pairs = [\n  [\"37\u00d790=3330\", \"45\u00d795=4275\"],\n  [\"47\u00d760=2820\", \"93\u00d734=3162\"],\n  [\"84\u00d790=7560\", \"67\u00d776=5092\"],\n  [\"35\u00d772=2520\", \"50\u00d750=2500\"],\n  [\"22\u00d789=1958\", \"19\u00d726=494\"],\n  [\"42\u00d736=1512\", \"78\u00d745=3510\"],\n  [\"13\u00d723=299\", \"32\u00d774=2368\"],\n  [\"21\u00d779=1659\", \"96\u00d753=5088\"],\n  [\"93\u00d791=8463\", \"61\u00d759=3599\"],\n  [\"74\u00d726=1924\", \"17\u00d760=1020\"],\n  [\"86\u00d791=7826\", \"55\u00d783=4565\"],\n  [\"61\u00d798=5978\", \"46\u00d722=1012\"],\n  [\"36\u00d797=3492\", \"63\u00d781=5103\"],\n  [\"19\u00d794=1786\", \"65\u00d794=6110\"],\n  [\"75\u00d775=5625\", \"45\u00d738=1710\"],\n  [\"89\u00d797=8633\", \"41\u00d793=3813\"],\n  [\"78\u00d739=3042\", \"39\u00d740=1560\"],\n  [\"74\u00d752=3848\", \"99\u00d725=2475\"],\n  [\"57\u00d716=912\", \"22\u00d786=1892\"],\n  [\"59\u00d755=3245\", \"64\u00d751=3264\"],\n  [\"48\u00d716=768\", \"54\u00d747=2538\"],\n  [\"19\u00d774=1406\", \"87\u00d739=3393\"],\n  [\"19\u00d771=1349\", \"90\u00d761=5490\"],\n  [\"87\u00d768=5916\", \"63\u00d782=5166\"],\n  [\"81\u00d742=3402\", \"32\u00d792=2944\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$d.Content.Find.Execute(\"37\u00d790=3330\", $false, $false, $false, $false, $false, $true, 1, $false, \"45\u00d795=4275\", 2) | Out-Null\n$d.Content.Find.Execute(\"47\u00d760=2820\", $false, $false, $false, $false, $false, $true, 1, $false, \"93\u00d734=3162\", 2) | Out-Null\n$d.Content.Find.Execute(\"84\u00d790=7560\", $false, $false, $false, $false, $false, $true, 1, $false, \"67\u00d776=5092\", 2) | Out-Null\n$d.Content.Find.Execute(\"35\u00d772=2520\", $false, $false, $false, $false, $false, $true, 1, $false, \"50\u00d750=2500\", 2) | Out-Null\n$d.Content.Find.Execute(\"22\u00d789=1958\", $false, $false, $false, $false, $false, $true, 1, $false, \"19\u00d726=494\", 2) | Out-Null\n$d.Content.Find.Execute(\"42\u00d736=1512\", $false, $false, $false, $false, $false, $true, 1, $false, \"78\u00d745=3510\", 2) | Out-Null\n$d.Content.Find.Execute(\"13\u00d723=299\", $false, $false, $false, $false, $false, $true, 1, $false, \"32\u00d774=2368\", 2) | Out-Null\n$d.Content.Find.Execute(\"21\u00d779=1659\", $false, $false, $false, $false, $false, $true, 1, $false, \"96\u00d753=5088\", 2) | Out-Null\n$d.Content.Find.Execute(\"93\u00d791=8463\", $false, $false, $false, $false, $false, $true, 1, $false, \"61\u00d759=3599\", 2) | Out-Null\n$d.Content.Find.Execute(\"74\u00d726=1924\", $false, $false, $false, $false, $false, $true, 1, $false, \"17\u00d760=1020\", 2) | Out-Null\n$d.Content.Find.Execute(\"86\u00d791=7826\", $false, $false, $false, $false, $false, $true, 1, $false, \"55\u00d783=4565\", 2) | Out-Null\n$d.Content.Find.Execute(\"61\u00d798=5978\", $false, $false, $false, $false, $false, $true, 1, $false, \"46\u00d722=1012\", 2) | Out-Null\n$d.Content.Find.Execute(\"36\u00d797=3492\", $false, $false, $false, $false, $false, $true, 1, $false, \"63\u00d781=5103\", 2) | Out-Null\n$d.Content.Find.Execute(\"19\u00d794=1786\", $false, $false, $false, $false, $false, $true, 1, $false, \"65\u00d794=6110\", 2) | Out-Null\n$d.Content.Find.Execute(\"75\u00d775=5625\", $false, $false, $false, $false, $false, $true, 1, $false, \"45\u00d738=1710\", 2) | Out-Null\n$d.Content.Find.Execute(\"89\u00d797=8633\", $false, $false, $false, $false, $false, $true, 1, $false, \"41\u00d793=3813\", 2) | Out-Null\n$d.Content.Find.Execute(\"78\u00d739=3042\", $false, $false, $false, $false, $false, $true, 1, $false, \"39\u00d740=1560\", 2) | Out-Null\n$d.Content.Find.Execute(\"74\u00d752=3848\", $false, $false, $false, $false, $false, $true, 1, $false, \"99\u00d725=2475\", 2) | Out-Null\n$d.Content.Find.Execute(\"57\u00d716=912\", $false, $false, $false, $false, $false, $true, 1, $false, \"22\u00d786=1892\", 2) | Out-Null\n$d.Content.Find.Execute(\"59\u00d755=3245\", $false, $false, $false, $false, $false, $true, 1, $false, \"64\u00d751=3264\", 2) | Out-Null\n$d.Content.Find.Execute(\"48\u00d716=768\", $false, $false, $false, $false, $false, $true, 1, $false, \"54\u00d747=2538\", 2) | Out-Null\n$d.Content.Find.Execute(\"19\u00d774=1406\", $false, $false, $false, $false, $false, $true, 1, $false, \"87\u00d739=3393\", 2) | Out-Null\n$d.Content.Find.Execute(\"19\u00d771=1349\", $false, $false, $false, $false, $false, $true, 1, $false, \"90\u00d761=5490\", 2) | Out-Null\n$d.Content.Find.Execute(\"87\u00d768=5916\", $false, $false, $false, $false, $false, $true, 1, $false, \"63\u00d782=5166\", 2) | Out-Null\n$d.Content.Find.Execute(\"81\u00d742=3402\", $false, $false, $false, $false, $false, $true, 1, $false, \"32\u00d792=2944\", 2) | Out-Null\n"}
